$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-07-19 Friday", $true, $false, $false, $false, $false, $false, 1, $false, "2024-07-20 Saturday", 2) | Out-Null
$d.Content.Find.Execute("38+33=", $true, $false, $false, $false, $false, $false, 1, $false, "95-56=", 2) | Out-Null
$d.Content.Find.Execute("41+13=", $true, $false, $false, $false, $false, $false, 1, $false, "75-50=", 2) | Out-Null
$d.Content.Find.Execute("98-21=", $true, $false, $false, $false, $false, $false, 1, $false, "18+61=", 2) | Out-Null
$d.Content.Find.Execute("77-8=", $true, $false, $false, $false, $false, $false, 1, $false, "31+10=", 2) | Out-Null
$d.Content.Find.Execute("29+70=", $true, $false, $false, $false, $false, $false, 1, $false, "6+27=", 2) | Out-Null
$d.Content.Find.Execute("37+54=", $true, $false, $false, $false, $false, $false, 1, $false, "92-42=", 2) | Out-Null
$d.Content.Find.Execute("50+11=", $true, $false, $false, $false, $false, $false, 1, $false, "72-0=", 2) | Out-Null
$d.Content.Find.Execute("1+17=", $true, $false, $false, $false, $false, $false, 1, $false, "58+38=", 2) | Out-Null
$d.Content.Find.Execute("53-27=", $true, $false, $false, $false, $false, $false, 1, $false, "98-1=", 2) | Out-Null
$d.Content.Find.Execute("6+92=", $true, $false, $false, $false, $false, $false, 1, $false, "19+71=", 2) | Out-Null
$d.Content.Find.Execute("26+43=", $true, $false, $false, $false, $false, $false, 1, $false, "37-19=", 2) | Out-Null
$d.Content.Find.Execute("16+50=", $true, $false, $false, $false, $false, $false, 1, $false, "56+7=", 2) | Out-Null
$d.Content.Find.Execute("60-38=", $true, $false, $false, $false, $false, $false, 1, $false, "89-32=", 2) | Out-Null
$d.Content.Find.Execute("21+77=", $true, $false, $false, $false, $false, $false, 1, $false, "16+26=", 2) | Out-Null
$d.Content.Find.Execute("3+84=", $true, $false, $false, $false, $false, $false, 1, $false, "43+17=", 2) | Out-Null
$d.Content.Find.Execute("89-60=", $true, $false, $false, $false, $false, $false, 1, $false, "65-49=", 2) | Out-Null
$d.Content.Find.Execute("87-37=", $true, $false, $false, $false, $false, $false, 1, $false, "72+3=", 2) | Out-Null
$d.Content.Find.Execute("74-72=", $true, $false, $false, $false, $false, $false, 1, $false, "47+4=", 2) | Out-Null
$d.Content.Find.Execute("93-51=", $true, $false, $false, $false, $false, $false, 1, $false, "59-38=", 2) | Out-Null
$d.Content.Find.Execute("39-10=", $true, $false, $false, $false, $false, $false, 1, $false, "24+70=", 2) | Out-Null
$d.Content.Find.Execute("22+29=", $true, $false, $false, $false, $false, $false, 1, $false, "63-10=", 2) | Out-Null
$d.Content.Find.Execute("42-19=", $true, $false, $false, $false, $false, $false, 1, $false, "36+55=", 2) | Out-Null
$d.Content.Find.Execute("79-8=", $true, $false, $false, $false, $false, $false, 1, $false, "15+82=", 2) | Out-Null
$d.Content.Find.Execute("82+1=", $true, $false, $false, $false, $false, $false, 1, $false, "91-7=", 2) | Out-Null
$d.Content.Find.Execute("28+38=", $true, $false, $false, $false, $false, $false, 1, $false, "38+10=", 2) | Out-Null
$d.Content.Find.Execute("11+80=", $true, $false, $false, $false, $false, $false, 1, $false, "7+26=", 2) | Out-Null
$d.Content.Find.Execute("94-21=", $true, $false, $false, $false, $false, $false, 1, $false, "28+8=", 2) | Out-Null
$d.Content.Find.Execute("30+54=", $true, $false, $false, $false, $false, $false, 1, $false, "74-4=", 2) | Out-Null
$d.Content.Find.Execute("54-2=", $true, $false, $false, $false, $false, $false, 1, $false, "73+6=", 2) | Out-Null
$d.Content.Find.Execute("75-35=", $true, $false, $false, $false, $false, $false, 1, $false, "88-41=", 2) | Out-Null
$d.Content.Find.Execute("35+35=", $true, $false, $false, $false, $false, $false, 1, $false, "51+19=", 2) | Out-Null
$d.Content.Find.Execute("50-23=", $true, $false, $false, $false, $false, $false, 1, $false, "47-45=", 2) | Out-Null
$d.Content.Find.Execute("8+80=", $true, $false, $false, $false, $false, $false, 1, $false, "87+8=", 2) | Out-Null
$d.Content.Find.Execute("21-14=", $true, $false, $false, $false, $false, $false, 1, $false, "82-77=", 2) | Out-Null
$d.Content.Find.Execute("34-6=", $true, $false, $false, $false, $false, $false, 1, $false, "81-14=", 2) | Out-Null
$d.Content.Find.Execute("4+34=", $true, $false, $false, $false, $false, $false, 1, $false, "65+30=", 2) | Out-Null
$d.Content.Find.Execute("7+64=", $true, $false, $false, $false, $false, $false, 1, $false, "38-1=", 2) | Out-Null
$d.Content.Find.Execute("72+8=", $true, $false, $false, $false, $false, $false, 1, $false, "93-25=", 2) | Out-Null
$d.Content.Find.Execute("54+42=", $true, $false, $false, $false, $false, $false, 1, $false, "66-22=", 2) | Out-Null
$d.Content.Find.Execute("68+17=", $true, $false, $false, $false, $false, $false, 1, $false, "99-15=", 2) | Out-Null
$d.Content.Find.Execute("91-32=", $true, $false, $false, $false, $false, $false, 1, $false, "48-31=", 2) | Out-Null
$d.Content.Find.Execute("55+41=", $true, $false, $false, $false, $false, $false, 1, $false, "95-38=", 2) | Out-Null
$d.Content.Find.Execute("29+2=", $true, $false, $false, $false, $false, $false, 1, $false, "44-15=", 2) | Out-Null
$d.Content.Find.Execute("47-10=", $true, $false, $false, $false, $false, $false, 1, $false, "94-93=", 2) | Out-Null
$d.Content.Find.Execute("66-60=", $true, $false, $false, $false, $false, $false, 1, $false, "15+33=", 2) | Out-Null
$d.Content.Find.Execute("82-12=", $true, $false, $false, $false, $false, $false, 1, $false, "46+3=", 2) | Out-Null
$d.Content.Find.Execute("96-93=", $true, $false, $false, $false, $false, $false, 1, $false, "26+37=", 2) | Out-Null
$d.Content.Find.Execute("67-23=", $true, $false, $false, $false, $false, $false, 1, $false, "59+1=", 2) | Out-Null
$d.Content.Find.Execute("31+11=", $true, $false, $false, $false, $false, $false, 1, $false, "74-2=", 2) | Out-Null
$d.Content.Find.Execute("94-42=", $true, $false, $false, $false, $false, $false, 1, $false, "74+9=", 2) | Out-Null
$d.Content.Find.Execute("46+32=", $true, $false, $false, $false, $false, $false, 1, $false, "16+59=", 2) | Out-Null
$d.Content.Find.Execute("21+68=", $true, $false, $false, $false, $false, $false, 1, $false, "29-24=", 2) | Out-Null
$d.Content.Find.Execute("59-8=", $true, $false, $false, $false, $false, $false, 1, $false, "18-5=", 2) | Out-Null
$d.Content.Find.Execute("92-38=", $true, $false, $false, $false, $false, $false, 1, $false, "41+45=", 2) | Out-Null
$d.Content.Find.Execute("53+28=", $true, $false, $false, $false, $false, $false, 1, $false, "9+59=", 2) | Out-Null
$d.Content.Find.Execute("67-48=", $true, $false, $false, $false, $false, $false, 1, $false, "98-70=", 2) | Out-Null
$d.Content.Find.Execute("10+22=", $true, $false, $false, $false, $false, $false, 1, $false, "19+50=", 2) | Out-Null
$d.Content.Find.Execute("48+4=", $true, $false, $false, $false, $false, $false, 1, $false, "8+66=", 2) | Out-Null
$d.Content.Find.Execute("74+16=", $true, $false, $false, $false, $false, $false, 1, $false, "51-18=", 2) | Out-Null
$d.Content.Find.Execute("83-54=", $true, $false, $false, $false, $false, $false, 1, $false, "85-41=", 2) | Out-Null
$d.Content.Find.Execute("13+19=", $true, $false, $false, $false, $false, $false, 1, $false, "88-88=", 2) | Out-Null
$d.Content.Find.Execute("87-82=", $true, $false, $false, $false, $false, $false, 1, $false, "0+40=", 2) | Out-Null
$d.Content.Find.Execute("65-12=", $true, $false, $false, $false, $false, $false, 1, $false, "16+8=", 2) | Out-Null
$d.Content.Find.Execute("14+35=", $true, $false, $false, $false, $false, $false, 1, $false, "18+33=", 2) | Out-Null
$d.Content.Find.Execute("55+5=", $true, $false, $false, $false, $false, $false, 1, $false, "78-20=", 2) | Out-Null
$d.Content.Find.Execute("79-46=", $true, $false, $false, $false, $false, $false, 1, $false, "0+0=", 2) | Out-Null
$d.Content.Find.Execute("71-0=", $true, $false, $false, $false, $false, $false, 1, $false, "77-48=", 2) | Out-Null
$d.Content.Find.Execute("80-62=", $true, $false, $false, $false, $false, $false, 1, $false, "2+89=", 2) | Out-Null
$d.Content.Find.Execute("35-25=", $true, $false, $false, $false, $false, $false, 1, $false, "72+1=", 2) | Out-Null
$d.Content.Find.Execute("89-13=", $true, $false, $false, $false, $false, $false, 1, $false, "82-14=", 2) | Out-Null
$d.Content.Find.Execute("53+30=", $true, $false, $false, $false, $false, $false, 1, $false, "2+41=", 2) | Out-Null
$d.Content.Find.Execute("51-28=", $true, $false, $false, $false, $false, $false, 1, $false, "63+8=", 2) | Out-Null
$d.Content.Find.Execute("40+30=", $true, $false, $false, $false, $false, $false, 1, $false, "49-5=", 2) | Out-Null
$d.Content.Find.Execute("70-57=", $true, $false, $false, $false, $false, $false, 1, $false, "58-21=", 2) | Out-Null
$d.Content.Find.Execute("13+39=", $true, $false, $false, $false, $false, $false, 1, $false, "22+32=", 2) | Out-Null
$d.Content.Find.Execute("27+54=", $true, $false, $false, $false, $false, $false, 1, $false, "21+64=", 2) | Out-Null
$d.Content.Find.Execute("26-21=", $true, $false, $false, $false, $false, $false, 1, $false, "73+0=", 2) | Out-Null
$d.Content.Find.Execute("23+71=", $true, $false, $false, $false, $false, $false, 1, $false, "8+51=", 2) | Out-Null
$d.Content.Find.Execute("16+1=", $true, $false, $false, $false, $false, $false, 1, $false, "84-18=", 2) | Out-Null
$d.Content.Find.Execute("42-41=", $true, $false, $false, $false, $false, $false, 1, $false, "41-14=", 2) | Out-Null
$d.Content.Find.Execute("36-7=", $true, $false, $false, $false, $false, $false, 1, $false, "95-67=", 2) | Out-Null
$d.Content.Find.Execute("27+30=", $true, $false, $false, $false, $false, $false, 1, $false, "51+30=", 2) | Out-Null
$d.Content.Find.Execute("93-30=", $true, $false, $false, $false, $false, $false, 1, $false, "17+69=", 2) | Out-Null
$d.Content.Find.Execute("88-26=", $true, $false, $false, $false, $false, $false, 1, $false, "15+69=", 2) | Out-Null
$d.Content.Find.Execute("93-29=", $true, $false, $false, $false, $false, $false, 1, $false, "25+58=", 2) | Out-Null
$d.Content.Find.Execute("86-11=", $true, $false, $false, $false, $false, $false, 1, $false, "51+26=", 2) | Out-Null
$d.Content.Find.Execute("82-56=", $true, $false, $false, $false, $false, $false, 1, $false, "19+73=", 2) | Out-Null
$d.Content.Find.Execute("22+6=", $true, $false, $false, $false, $false, $false, 1, $false, "48+1=", 2) | Out-Null
$d.Content.Find.Execute("94+3=", $true, $false, $false, $false, $false, $false, 1, $false, "51+43=", 2) | Out-Null
$d.Content.Find.Execute("95-11=", $true, $false, $false, $false, $false, $false, 1, $false, "63-14=", 2) | Out-Null
$d.Content.Find.Execute("20+5=", $true, $false, $false, $false, $false, $false, 1, $false, "35-20=", 2) | Out-Null
$d.Content.Find.Execute("18+22=", $true, $false, $false, $false, $false, $false, 1, $false, "88-5=", 2) | Out-Null
$d.Content.Find.Execute("76-63=", $true, $false, $false, $false, $false, $false, 1, $false, "39+51=", 2) | Out-Null
$d.Content.Find.Execute("47-38=", $true, $false, $false, $false, $false, $false, 1, $false, "51-46=", 2) | Out-Null
$d.Content.Find.Execute("59+27=", $true, $false, $false, $false, $false, $false, 1, $false, "0+71=", 2) | Out-Null
$d.Content.Find.Execute("15+45=", $true, $false, $false, $false, $false, $false, 1, $false, "68-45=", 2) | Out-Null
$d.Content.Find.Execute("33+48=", $true, $false, $false, $false, $false, $false, 1, $false, "70-53=", 2) | Out-Null
$d.Content.Find.Execute("20+44=", $true, $false, $false, $false, $false, $false, 1, $false, "54-47=", 2) | Out-Null
$d.Content.Find.Execute("99-55=", $true, $false, $false, $false, $false, $false, 1, $false, "13+69=", 2) | Out-Null
$d.Content.Find.Execute("35-23=", $true, $false, $false, $false, $false, $false, 1, $false, "40+59=", 2) | Out-Null
